$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.489.18"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = "'  +0.92%  "
$ws.Range("E2").Style = 'Normal'
$ws.Range("D3").Value = "'3.388.81"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = "'  -0.01%  "
$ws.Range("E3").Style = 'Normal'
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = 'Normal'
$ws.Range("D5").Value = "'575.31"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = "'  +0.57%  "
$ws.Range("E5").Style = 'Normal'
$ws.Range("D6").Value = "'140.85"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = "'  -0.95%  "
$ws.Range("E6").Style = 'Normal'
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = 'Normal'
$ws.Range("D8").Value = "'0.474"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = "'  -0.39%  "
$ws.Range("E8").Style = 'Normal'
$ws.Range("D9").Value = "'7.72"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = "'  +2.35%  "
$ws.Range("E9").Style = 'Normal'
$ws.Range("E10").Value = "'  -0.90%  "
$ws.Range("E10").Style = 'Normal'
$ws.Range("D11").Value = "'0.387"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = "'  -2.19%  "
$ws.Range("E11").Style = 'Normal'
$ws.Range("D12").Value = "'3.972.70"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = "'  +0.17%  "
$ws.Range("E12").Style = 'Normal'
$ws.Range("E13").Value = "'  +0.03%  "
$ws.Range("E13").Style = 'Normal'
$ws.Range("D14").Value = "'28.37"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = "'  +0.63%  "
$ws.Range("E14").Style = 'Normal'
$ws.Range("D15").Value = "'3.389.09"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = "'  +0.00%  "
$ws.Range("E15").Style = 'Normal'
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = "'  -0.32%  "
$ws.Range("E16").Style = 'Normal'
$ws.Range("D17").Value = "'61.536.64"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = "'  +0.92%  "
$ws.Range("E17").Style = 'Normal'
$ws.Range("D18").Value = "'6.13"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = "'  -0.55%  "
$ws.Range("E18").Style = 'Normal'
$ws.Range("D19").Value = "'13.65"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = "'  -1.24%  "
$ws.Range("E19").Style = 'Normal'
$ws.Range("D20").Value = "'8.97"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = "'  +0.00%  "
$ws.Range("E20").Style = 'Normal'
$ws.Range("D21").Value = "'390.41"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = "'  +1.65%  "
$ws.Range("E21").Style = 'Normal'
$ws.Range("D22").Value = "'75.27"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E23").Value = "'  -0.78%  "
$ws.Range("E23").Style = 'Normal'
$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("E24").Style = 'Normal'
$ws.Range("D25").Value = "'0.0000112"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = "'  -4.30%  "
$ws.Range("E25").Style = 'Normal'
$ws.Range("D26").Value = "'0.193"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = "'  +7.57%  "
$ws.Range("E26").Style = 'Normal'
$ws.Range("E27").Value = "'  +0.01%  "
$ws.Range("E27").Style = 'Normal'
$ws.Range("E28").Value = "'  -1.76%  "
$ws.Range("E28").Style = 'Normal'
$ws.Range("E29").Value = "'  +0.37%  "
$ws.Range("E29").Style = 'Normal'
$ws.Range("E30").Value = "'  +0.01%  "
$ws.Range("E30").Style = 'Normal'
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("E31").Style = 'Normal'
$ws.Range("E32").Value = "'  -3.24%  "
$ws.Range("E32").Style = 'Normal'
$ws.Range("D33").Value = "'23.29"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = "'  -0.96%  "
$ws.Range("E33").Style = 'Normal'
$ws.Range("D34").Value = "'6.92"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = "'  -0.98%  "
$ws.Range("E34").Style = 'Normal'
$ws.Range("D35").Value = "'168.09"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = "'  +0.18%  "
$ws.Range("E35").Style = 'Normal'
$ws.Range("D36").Value = "'5.05"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = "'  +1.25%  "
$ws.Range("E36").Style = 'Normal'
$ws.Range("D37").Value = "'3.427.65"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = "'  +0.31%  "
$ws.Range("E37").Style = 'Normal'
$ws.Range("E38").Value = "'  -1.40%  "
$ws.Range("E38").Style = 'Normal'
$ws.Range("D39").Value = "'0.0771"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = "'  -0.53%  "
$ws.Range("E39").Style = 'Normal'
$ws.Range("E40").Value = "'  -5.70%  "
$ws.Range("E40").Style = 'Normal'
$ws.Range("D41").Value = "'0.779"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = "'  -0.12%  "
$ws.Range("E41").Style = 'Normal'
$ws.Range("E42").Value = "'  -0.14%  "
$ws.Range("E42").Style = 'Normal'
$ws.Range("E43").Value = "'  -1.01%  "
$ws.Range("E43").Style = 'Normal'
$ws.Range("E44").Value = "'  +0.76%  "
$ws.Range("E44").Style = 'Normal'
$ws.Range("D45").Value = "'2.468.33"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = "'  -0.40%  "
$ws.Range("E45").Style = 'Normal'
$ws.Range("D46").Value = "'22.96"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = "'  -0.18%  "
$ws.Range("E46").Style = 'Normal'
$ws.Range("D47").Value = "'6.66"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = "'  -2.27%  "
$ws.Range("E47").Style = 'Normal'
$ws.Range("E48").Value = "'  +0.07%  "
$ws.Range("E48").Style = 'Normal'
$ws.Range("E49").Value = "'  -1.70%  "
$ws.Range("E49").Style = 'Normal'
$ws.Range("D50").Value = "'2.05"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = "'  -1.95%  "
$ws.Range("E50").Style = 'Normal'
$ws.Range("D51").Value = "'0.205"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = "'  -1.93%  "
$ws.Range("E51").Style = 'Normal'
